$wb = $excel.ActiveWorkbook

# --- 1. Metadata sheet: bump "Last Updated" timestamp by one minute ---
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("A2").Value = "05 Nov 2025, 10:21 AM"

# --- 2. Industry Analysis sheet: update "1 Year" (column F) figures ---
$wsIndustry = $wb.Worksheets.Item("Industry Analysis")
$industryUpdates = @(
    @(2, 21.3),
    @(3, -4.3927),
    @(4, 35.9445),
    @(5, -51.0482),
    @(6, 57.2275),
    @(7, -9.640700000000001),
    @(8, -6.1449),
    @(9, 36.9733),
    @(10, -4.7026),
    @(11, 46.5317),
    @(12, -2.102),
    @(13, 17.4681),
    @(14, -33.0245),
    @(15, 1.0205),
    @(16, 2.0426),
    @(17, -16.2411),
    @(18, 7.4627),
    @(19, -25.798),
    @(20, 47.7485),
    @(21, 19.5587),
    @(22, 76.5603),
    @(23, -54.2675),
    @(24, -0.8811),
    @(25, 4.8518),
    @(26, 3.6831),
    @(27, -34.0874),
    @(28, -11.9893),
    @(29, -12.994),
    @(30, 25.5415),
    @(31, 56.5088),
    @(32, 2.0908),
    @(33, -4.7193),
    @(34, 22.8807),
    @(35, 5.3359),
    @(36, -5.1995),
    @(37, -5.6238),
    @(38, -22.595),
    @(39, 10.8405),
    @(40, -7.5963),
    @(41, -4.552),
    @(42, 22.3098),
    @(43, 14.0694),
    @(44, -9.6066),
    @(45, 27.639),
    @(46, -6.3484),
    @(47, -40.5302),
    @(48, -29.7988),
    @(49, -24.0791),
    @(50, -49.1803),
    @(51, -51.6023),
    @(52, -34.4756),
    @(53, -11.5478),
    @(54, -2.3796),
    @(55, -15.4382),
    @(56, -27.6987),
    @(57, -27.1559),
    @(58, -2.1585),
    @(59, -23.0964),
    @(60, -13.3217),
    @(61, -8.1496),
    @(62, -16.0695),
    @(63, -12.5465),
    @(64, 47.7264),
    @(65, -42.4232),
    @(66, 11.3291),
    @(67, 14.3746),
    @(68, 32.6702),
    @(69, -17.0097),
    @(70, -13.5162),
    @(71, 11.4259),
    @(72, 2.6754),
    @(73, -11.1574),
    @(74, -13.2502),
    @(75, 24.7078),
    @(76, 53.3554)
)
foreach ($row in $industryUpdates) {
    $r = $row[0]
    $val = $row[1]
    $wsIndustry.Cells.Item($r, 6).Value = $val
}

# --- 3. Stock List sheet: refresh the stock roster (rows 2-76) ---
$wsStockList = $wb.Worksheets.Item("Stock List")
$stockListUpdates = @(
    @(2, "NIFTYCASE", 10.19, -0.5854, 0),
    @(3, "MOMENTUM30", 31.54, -0.6614, 0),
    @(4, "CANHLIFE", 118.46, 0.6286, 11253.7),
    @(5, "FLEXIADD", 10.64, -1.0233, 0),
    @(6, "MOENERGY", 36.3, -0.6568000000000001, 0),
    @(7, "MONIFTY100", 26.49, 0.3409, 0),
    @(8, "RUBICON", 652.65, -0.1453, 10752.4289),
    @(9, "CRAMC", 317.2, 2.3226, 6325.5208),
    @(10, "LGEINDIA", 1633.4, -0.946, 110870.6825),
    @(11, "TATACAP", 329.3, 0.1521, 139783.5374),
    @(12, "ELIQUID", 1004.85, 0.0408, 0),
    @(13, "WEWORK", 632.15, -2.4008, 8472.2803),
    @(14, "GROWWRLTY", 10.8, -0.4608, 0),
    @(15, "ADVANCE", 130.05, -5.2666, 836.0358),
    @(16, "OMFREIGHT", 88.90000000000001, -0.5926, 299.3747),
    @(17, "GLOTTIS", 72.73999999999999, -0.8587, 672.1394),
    @(18, "FABTECH", 237.72, 0.4734, 1056.6843),
    @(19, "PACEDIGITK", 218.85, 0.1327, 4723.9063),
    @(20, "JAINREC", 377.25, 1.2208, 13018.3623),
    @(21, "EPACKPEB", 301.45, 1.979, 3028.1254),
    @(22, "BMWVENTLTD", 69.25, 0, 600.5014),
    @(23, "STYL", 372.4, -0.8388, 6025.649),
    @(24, "JARO", 621.5, -1.4821, 1377.0134),
    @(25, "SOLARWORLD", 309.1, -0.6269, 2679.0517),
    @(26, "ARSSBL", 537.3, 4.7266, 3370.2277),
    @(27, "GANESHCP", 274.4, -2.7984, 1108.9312),
    @(28, "ATLANTAELE", 1003.05, -1.7436, 7713.116),
    @(29, "GKENERGY", 213.85, -0.7933, 4337.2472),
    @(30, "SAATVIKGL", 528.2, -1.3079, 6713.6863),
    @(31, "IVALUE", 281.45, -0.3364, 1506.8799),
    @(32, "VMSTMT", 70.03, -0.9056, 347.5674),
    @(33, "EUROPRATIK", 321.75, 0.8147, 3288.285),
    @(34, "SHRINGARMS", 229.31, -1.2616, 2211.284),
    @(35, "DEVX", 44.53, -0.3803, 401.605),
    @(36, "URBANCO", 148.9, -2.0459, 21380.5798),
    @(37, "SML100CASE", 10.36, -0.7663, 0),
    @(38, "AONEGOLD", 11.28, -0.2653, 0),
    @(39, "ELM250", 16.72, 0.1797, 0),
    @(40, "AMANTA", 122.52, 1.407, 475.7372),
    @(41, "CPEDU", 315.9, 1.8539, 574.7148999999999),
    @(42, "AHCL", 139.27, 3.1706, 740.2409),
    @(43, "STLNETWORK", 26.59, -0.412, 1297.3822),
    @(44, "VIKRAN", 98.05, -1.783, 2528.8166),
    @(45, "MANUFGBEES", 151.77, -1.011, 0),
    @(46, "MEIL", 461.15, -0.7319, 1274.1632),
    @(47, "GROWWNXT50", 70.29000000000001, -0.4109, 0),
    @(48, "SHREEJISPG", 270.05, -0.7899, 4399.6074),
    @(49, "GEMAROMA", 219.52, -0.876, 1146.7097),
    @(50, "PATELRMART", 219.31, -1.0646, 732.5069999999999),
    @(51, "VIKRAMSOLR", 322, -1.5892, 11647.2884),
    @(52, "LTGILTCASE", 29.67, 0.2365, 0),
    @(53, "REGAAL", 89.13, -0.8675, 915.5742),
    @(54, "BLUESTONE", 711.95, 0.1266, 10773.2539),
    @(55, "MOSILVER", 145.9, -1.5054, 0),
    @(56, "ALLTIME", 308.75, 2.66, 2022.5526),
    @(57, "JSWCEMENT", 134.98, -0.4793, 18402.6999),
    @(58, "SBILIQETF", 1012.94, 0.0296, 0),
    @(59, "HILINFRA", 77.23, -0.3998, 0),
    @(60, "GROWWPOWER", 10.28, -0.9634, 0),
    @(61, "LOTUSDEV", 177.82, 0.3669, 8690.485000000001),
    @(62, "MBEL", 450.2, -0.7714, 2572.8126),
    @(63, "LAXMIINDIA", 145.62, -1.1942, 761.1248000000001),
    @(64, "CPPLUS", 1322.1, -0.264, 15497.9053),
    @(65, "SHANTIGOLD", 241.57, -1.6409, 1741.6231),
    @(66, "MOGOLD", 119.65, -0.5403, 0),
    @(67, "BRIGHOTEL", 82.39, -0.9855, 3129.5229),
    @(68, "INDIQUBE", 212.64, -0.7561, 4465.6847),
    @(69, "EBGNG", 346.65, 3.2311, 3952.2092),
    @(70, "LIQGRWBEES", 1014.74, 0.0246, 0),
    @(71, "CHEMBONDCH", 153.35, -1.6987, 412.459),
    @(72, "GROWWNIFTY", 10.29, -0.3872, 0),
    @(73, "ANTHEM", 702.25, -0.1209, 39439.0658),
    @(74, "QUALITY30", 21.05, -0.8945, 0),
    @(75, "SMARTWORKS", 606.65, 2.0867, 6931.2448),
    @(76, "TRAVELFOOD", 1316.3, 0.1141, 17332.9705)
)
foreach ($row in $stockListUpdates) {
    $r = $row[0]
    $ticker = $row[1]
    $price = $row[2]
    $pctChange = $row[3]
    $marketCap = $row[4]
    $wsStockList.Cells.Item($r, 2).Value = $ticker
    $wsStockList.Cells.Item($r, 3).Value = $ticker
    $wsStockList.Cells.Item($r, 4).Value = $price
    $wsStockList.Cells.Item($r, 5).Value = $pctChange
    $wsStockList.Cells.Item($r, 8).Value = $marketCap
}
